$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.691.32"
$ws.Range("E2").Value = "  +3.63%  "
$ws.Range("D3").Value = "2.422.13"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'315.93"
$ws.Range("E5").Value = "  +3.89%  "
$ws.Range("D6").Value = "'101.70"
$ws.Range("E6").Value = "  +6.47%  "
$ws.Range("D7").Value = "'0.515"
$ws.Range("E7").Value = "  +2.47%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +9.62%  "
$ws.Range("D10").Value = "'35.45"
$ws.Range("E10").Value = "  +3.14%  "
$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = "  +1.77%  "
$ws.Range("D12").Value = "'18.79"
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  -1.59%  "
$ws.Range("D14").Value = "'6.94"
$ws.Range("E14").Value = "  +3.19%  "
$ws.Range("D15").Value = "2.804.65"
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").Value = "2.439.92"
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "'0.832"
$ws.Range("E17").Value = "  +4.20%  "
$ws.Range("D18").Value = "44.574.06"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("D19").Value = "'12.29"
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").Value = "'6.36"
$ws.Range("E20").Value = "  +1.03%  "
$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("E21").Value = "  +3.53%  "
$ws.Range("D22").Value = "'68.73"
$ws.Range("D23").Value = "'242.43"
$ws.Range("E23").Value = "  +2.94%  "
$ws.Range("D24").Value = "'2.28"
$ws.Range("E24").Value = "  +3.90%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("E25").Value = "  +1.58%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +2.44%  "
$ws.Range("D28").Value = "'2.18"
$ws.Range("E28").Value = "  -8.06%  "
$ws.Range("D29").Value = "'9.53"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").Value = "'33.65"
$ws.Range("E30").Value = "  +3.57%  "
$ws.Range("D31").Value = "'48.54"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "'0.126"
$ws.Range("E32").Value = "  +17.59%  "
$ws.Range("D33").Value = "'19.50"
$ws.Range("E33").Value = "  +11.04%  "
$ws.Range("D34").Value = "'5.16"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'0.0775"
$ws.Range("E35").Value = "  +6.28%  "
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  +2.66%  "
$ws.Range("D38").Value = "'4.48"
$ws.Range("E38").Value = "  +3.03%  "
$ws.Range("D39").Value = "'2.87"
$ws.Range("E39").Value = "  +0.78%  "
$ws.Range("D40").Value = "'123.78"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'0.109"
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -2.92%  "
$ws.Range("D43").Value = "'20.86"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("E44").Value = "  +4.39%  "
$ws.Range("D45").Value = "1.943.19"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  +8.02%  "
$ws.Range("D48").Value = "'9.48"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").Value = "'1.72"
$ws.Range("E49").Value = "  +13.55%  "
$ws.Range("D50").Value = "'75.18"
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("D51").Value = "'54.01"
$ws.Range("E51").Value = "  +5.24%  "
